$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set row heights for the new rows
$ws.Rows.Item(57).RowHeight = 187
$ws.Rows.Item(58).RowHeight = 68
$ws.Rows.Item(59).RowHeight = 42
$ws.Rows.Item(60).RowHeight = 42

# B57 -> si index 142
$ws.Range("B49").Copy()
$ws.Range("B57").PasteSpecial(-4122)
$c = $ws.Range("B57")
$c.Value2 = '拷贝数组的方式'

# C57 -> si index 151
$ws.Range("C3").Copy()
$ws.Range("C57").PasteSpecial(-4122)
$c = $ws.Range("C57")
$c.Value2 = '浅拷贝' + "`n" + 'slice 方法返回一个新的数组对象，这一对象是一个由 begin 和 end 决定的原数组的浅拷贝（包括 begin，不包括end）。如果该元素是个对象引用 （不是实际的对象），slice 会拷贝这个对象引用到新的数组里。两个对象引用都引用了同一个对象。如果被引用的对象发生改变，则新的和原来的数组中的这个元素也会发生改变。元素的其他数据类型正常拷贝'
$c.Characters(1,3).Font.Bold = $true
$c.Characters(52,3).Font.Bold = $true
$c.Characters(80,14).Font.Bold = $true

# D57 -> si index 152
$ws.Range("C3").Copy()
$ws.Range("D57").PasteSpecial(-4122)
$ws.Range("D57").VerticalAlignment = -4160
$c = $ws.Range("D57")
$c.Value2 = '浅拷贝' + "`n" + 'concat 方法用于合并两个或多个数组。此方法不会更改现有数组，而是返回一个新数组. 对象引用（而不是实际对象）：concat将对象引用复制到新数组中。 原始数组和新数组都引用相同的对象。 也就是说，如果引用的对象被修改，则更改对于新数组和原始数组都是可见的。 这包括也是数组的数组参数的元素。元素的其他数据类型正常拷贝.'
$c.Characters(1,3).Font.Bold = $true
$c.Characters(49,13).Font.Bold = $true

# E57 -> si index 143
$ws.Range("C3").Copy()
$ws.Range("E57").PasteSpecial(-4122)
$ws.Range("E57").VerticalAlignment = -4160
$c = $ws.Range("E57")
$c.Value2 = '浅拷贝' + "`n" + 'es6 展开运算符 (...) 可以在函数调用/数组构造时, 将数组表达式或者string在语法层面展开；还可以在构造字面量对象时, 将对象表达式按key-value的方式展开' + "`n" + ''

# F57 -> si index 144
$ws.Range("C3").Copy()
$ws.Range("F57").PasteSpecial(-4122)
$ws.Range("F57").VerticalAlignment = -4160
$c = $ws.Range("F57")
$c.Value2 = '深拷贝' + "`n" + 'json.parse(json.stringify(arr))'

# B58 -> si index 149
$ws.Range("B49").Copy()
$ws.Range("B58").PasteSpecial(-4122)
$c = $ws.Range("B58")
$c.Value2 = '数组遍历的几张方法返回值'

# C58 -> si index 150
$ws.Range("C3").Copy()
$ws.Range("C58").PasteSpecial(-4122)
$c = $ws.Range("C58")
$c.Value2 = 'map 返回各元素是否满足条件的boolean 数组[true, false, undefined]'

# D58 -> si index 145
$ws.Range("C3").Copy()
$ws.Range("D58").PasteSpecial(-4122)
$c = $ws.Range("D58")
$c.Value2 = 'filter 返回符合过滤条件的部分数组'

# E58 -> si index 146
$ws.Range("C3").Copy()
$ws.Range("E58").PasteSpecial(-4122)
$c = $ws.Range("E58")
$c.Value2 = 'some 只要有任意一个元素满足条件，返回 true ,否则返回 false'

# F58 -> si index 147
$ws.Range("C3").Copy()
$ws.Range("F58").PasteSpecial(-4122)
$c = $ws.Range("F58")
$c.Value2 = 'every 必须全部元素满足条件，返回true，否则返回 false'

# G58 -> si index 148
$ws.Range("C3").Copy()
$ws.Range("G58").PasteSpecial(-4122)
$c = $ws.Range("G58")
$c.Value2 = 'forEach 总是返回undefined'

# B59 -> si index 153
$ws.Range("B49").Copy()
$ws.Range("B59").PasteSpecial(-4122)
$c = $ws.Range("B59")
$c.Value2 = 'null 和 undefined区别'

# C59 -> si index 156
$ws.Range("C3").Copy()
$ws.Range("C59").PasteSpecial(-4122)
$c = $ws.Range("C59")
$c.Value2 = 'null 表示一个值被定义了，但这个值是空值 false'

# D59 -> si index 154
$ws.Range("C3").Copy()
$ws.Range("D59").PasteSpecial(-4122)
$c = $ws.Range("D59")
$c.Value2 = 'undefined 表示此处应该有值，但还没定义, false'

# E59 -> si index 155
$ws.Range("C3").Copy()
$ws.Range("E59").PasteSpecial(-4122)
$c = $ws.Range("E59")
$c.Value2 = 'null == undefined  结果 true' + "`n" + 'null === undefined 结果 false'

# B60 -> si index 157
$ws.Range("B49").Copy()
$ws.Range("B60").PasteSpecial(-4122)
$c = $ws.Range("B60")
$c.Value2 = 'vue 组件中name的用处'

# C60 -> si index 158
$ws.Range("C3").Copy()
$ws.Range("C60").PasteSpecial(-4122)
$c = $ws.Range("C60")
$c.Value2 = '1. keep-alive 时用来标记组件'

# D60 -> si index 159
$ws.Range("C3").Copy()
$ws.Range("D60").PasteSpecial(-4122)
$c = $ws.Range("D60")
$c.Value2 = '2. 递归组件时用'

# E60 -> si index 160
$ws.Range("C3").Copy()
$ws.Range("E60").PasteSpecial(-4122)
$c = $ws.Range("E60")
$c.Value2 = '3. vue 浏览器插件调试时查看组件名'
